$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.635.62'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '2.285.18'
$ws.Range('E3').Value = '  +1.91%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '95.48'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '266.74'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -1.72%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.608'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -5.75%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.17'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -8.73%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0935'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.39%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.73'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -7.98%  '
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '2.618.01'
$ws.Range('E14').Value = '  +3.20%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.17'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.849'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +3.52%  '
$ws.Range('D17').Value = '2.293.49'
$ws.Range('E17').Value = '  +1.88%  '
$ws.Range('D18').Value = '43.602.07'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000107'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.18'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '72.29'
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.35'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.80%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.41'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.92'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -3.44%  '
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.30'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.72%  '
$ws.Range('E27').Value = '  -1.33%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '3.45'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('E29').Value = '  +0.35%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '175.66'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '37.84'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -7.13%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '21.95'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +4.30%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0882'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.35%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.36'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.45%  '
$ws.Range('E35').Value = '  +0.57%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.108'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.77%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0353'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.41'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.87%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.25'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -12.90%  '
$ws.Range('E40').Value = '  +7.86%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.237'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -6.16%  '
$ws.Range('E42').Value = '  +16.23%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '11.85'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -8.21%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '62.18'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.14%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.83'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +4.59%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.21'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '98.11'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.66%  '
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.55'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +8.38%  '
$ws.Range('D51').Value = '2.507.07'
$ws.Range('E51').Value = '  +2.24%  '
